$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '75.937.11'
$ws.Range("E2").Value = '  +1.43%  '
$ws.Range("D3").Value = '2.930.68'
$ws.Range("E3").Value = '  +3.98%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '203.33'
$ws.Range("E5").Value = '  +8.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '596.44'
$ws.Range("E6").Value = '  +0.69%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +0.52%  '
$ws.Range("E9").Value = '  +3.84%  '
$ws.Range("D10").Value = '2.932.96'
$ws.Range("E10").Value = '  +4.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.435'
$ws.Range("E12").Value = '  +0.72%  '
$ws.Range("E13").Value = '  +0.94%  '
$ws.Range("D14").Value = '3.474.79'
$ws.Range("E14").Value = '  +4.27%  '
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '75.974.38'
$ws.Range("E15").Value = '  +1.48%  '
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.06'
$ws.Range("E16").Value = '  +4.70%  '
$ws.Range("E17").Value = '  +2.05%  '
$ws.Range("D18").Value = '2.940.46'
$ws.Range("E18").Value = '  +4.49%  '
$ws.Range("E19").Value = '  +7.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.87'
$ws.Range("E20").Value = '  -1.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '372.55'
$ws.Range("E21").Value = '  -1.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.31'
$ws.Range("E22").Value = '  +2.04%  '
$ws.Range("E23").Value = '  +5.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.76'
$ws.Range("E24").Value = '  +1.24%  '
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("D26").Value = '3.087.07'
$ws.Range("E26").Value = '  +4.75%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.31'
$ws.Range("E27").Value = '  +3.94%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.68'
$ws.Range("E28").Value = '  -0.51%  '
$ws.Range("E29").Value = '  +4.86%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.39'
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '502.79'
$ws.Range("E32").Value = '  -2.08%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.80'
$ws.Range("E33").Value = '  +1.79%  '
$ws.Range("E34").Value = '  +2.90%  '
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '20.29'
$ws.Range("E36").Value = '  +1.87%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.85'
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("B38").Value = 'Cronos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.111'
$ws.Range("E38").Value = '  +28.81%  '
$ws.Range("E39").Value = '  +1.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.373'
$ws.Range("E40").Value = '  +9.34%  '
$ws.Range("E41").Value = '  -4.71%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '181.73'
$ws.Range("E42").Value = '  -2.18%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("E44").Value = '  +0.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.66'
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.13'
$ws.Range("E46").Value = '  +0.29%  '
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("E48").Value = '  +1.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.580'
$ws.Range("E49").Value = '  +0.77%  '
$ws.Range("E50").Value = '  +0.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.50'
$ws.Range("E51").Value = '  +8.06%  '
